$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column header in H1 — copy the existing header formatting
# (bold/border/centered, same style as the other header cells such as G1)
# onto H1 before writing its text, then set the value.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("H1").Value = "Save"

# Fill the new "Save" column with data
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
